$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Update Hoja1 (sheet1) values, columns B, C, D, rows 2-27 ---
$ws1.Range("B2").Value = 14545
$ws1.Range("C2").Value = 13258
$ws1.Range("D2").Value = 1472
$ws1.Range("B3").Value = 13040
$ws1.Range("C3").Value = 14101
$ws1.Range("D3").Value = 525
$ws1.Range("B4").Value = 3722
$ws1.Range("C4").Value = 3078
$ws1.Range("D4").Value = 536
$ws1.Range("B5").Value = 6649
$ws1.Range("C5").Value = 6894
$ws1.Range("D5").Value = 3709
$ws1.Range("B6").Value = 3739
$ws1.Range("C6").Value = 4357
$ws1.Range("D6").Value = 948
$ws1.Range("B7").Value = 5869
$ws1.Range("C7").Value = 7050
$ws1.Range("D7").Value = 462
$ws1.Range("B8").Value = 6292
$ws1.Range("C8").Value = 6903
$ws1.Range("D8").Value = 1639
$ws1.Range("B9").Value = 3510
$ws1.Range("C9").Value = 3979
$ws1.Range("D9").Value = 722
$ws1.Range("B10").Value = 24797
$ws1.Range("C10").Value = 27500
$ws1.Range("D10").Value = 5640
$ws1.Range("B11").Value = 2269
$ws1.Range("C11").Value = 2789
$ws1.Range("D11").Value = 400
$ws1.Range("B12").Value = 2540
$ws1.Range("C12").Value = 2192
$ws1.Range("D12").Value = 311
$ws1.Range("B13").Value = 87103
$ws1.Range("C13").Value = 87168
$ws1.Range("D13").Value = 59161
$ws1.Range("B14").Value = 35922
$ws1.Range("C14").Value = 38453
$ws1.Range("D14").Value = 7396
$ws1.Range("B15").Value = 1254
$ws1.Range("C15").Value = 1668
$ws1.Range("D15").Value = 213
$ws1.Range("B16").Value = 8300
$ws1.Range("C16").Value = 8423
$ws1.Range("D16").Value = 869
$ws1.Range("B17").Value = 15021
$ws1.Range("C17").Value = 16187
$ws1.Range("D17").Value = 3904
$ws1.Range("B18").Value = 12150
$ws1.Range("C18").Value = 21091
$ws1.Range("D18").Value = 3420
$ws1.Range("B19").Value = 405
$ws1.Range("C19").Value = 587
$ws1.Range("D19").Value = 224
$ws1.Range("B20").Value = 7402
$ws1.Range("C20").Value = 8490
$ws1.Range("D20").Value = 922
$ws1.Range("B21").Value = 10048
$ws1.Range("C21").Value = 12134
$ws1.Range("D21").Value = 5634
$ws1.Range("B22").Value = 4487
$ws1.Range("C22").Value = 4820
$ws1.Range("D22").Value = 449
$ws1.Range("B23").Value = 7312
$ws1.Range("C23").Value = 7599
$ws1.Range("D23").Value = 3982
$ws1.Range("B24").Value = 5676
$ws1.Range("C24").Value = 5559
$ws1.Range("D24").Value = 1495
$ws1.Range("B25").Value = 307748
$ws1.Range("C25").Value = 319321
$ws1.Range("D25").Value = 181544
$ws1.Range("B26").Value = 9012
$ws1.Range("C26").Value = 9208
$ws1.Range("D26").Value = 1155
$ws1.Range("B27").Value = 50556
$ws1.Range("C27").Value = 61673
$ws1.Range("D27").Value = 3753

# --- Update Hoja2 (sheet2) values, columns A, B, C, rows 2-19 ---
$ws2.Range("A2").Value = 19146
$ws2.Range("B2").Value = 19419
$ws2.Range("C2").Value = 8546
$ws2.Range("A3").Value = 17870
$ws2.Range("B3").Value = 19304
$ws2.Range("C3").Value = 10871
$ws2.Range("A4").Value = 8073
$ws2.Range("B4").Value = 8426
$ws2.Range("C4").Value = 3856
$ws2.Range("A5").Value = 4062
$ws2.Range("B5").Value = 4652
$ws2.Range("C5").Value = 1545
$ws2.Range("A6").Value = 4038
$ws2.Range("B6").Value = 4709
$ws2.Range("C6").Value = 1346
$ws2.Range("A7").Value = 5275
$ws2.Range("B7").Value = 6169
$ws2.Range("C7").Value = 2971
$ws2.Range("A8").Value = 7982
$ws2.Range("B8").Value = 8276
$ws2.Range("C8").Value = 4092
$ws2.Range("A9").Value = 18075
$ws2.Range("B9").Value = 19205
$ws2.Range("C9").Value = 10109
$ws2.Range("A10").Value = 53272
$ws2.Range("B10").Value = 54836
$ws2.Range("C10").Value = 33579
$ws2.Range("A11").Value = 16887
$ws2.Range("B11").Value = 16777
$ws2.Range("C11").Value = 9246
$ws2.Range("A12").Value = 3371
$ws2.Range("B12").Value = 3287
$ws2.Range("C12").Value = 1563
$ws2.Range("A13").Value = 10138
$ws2.Range("B13").Value = 10322
$ws2.Range("C13").Value = 4312
$ws2.Range("A14").Value = 79656
$ws2.Range("B14").Value = 82236
$ws2.Range("C14").Value = 43842
$ws2.Range("A15").Value = 3286
$ws2.Range("B15").Value = 2854
$ws2.Range("C15").Value = 1190
$ws2.Range("A16").Value = 10761
$ws2.Range("B16").Value = 11134
$ws2.Range("C16").Value = 8084
$ws2.Range("A17").Value = 40495
$ws2.Range("B17").Value = 42360
$ws2.Range("C17").Value = 32523
$ws2.Range("A18").Value = 4645
$ws2.Range("B18").Value = 4863
$ws2.Range("C18").Value = 3556
$ws2.Range("A19").Value = 716
$ws2.Range("B19").Value = 492
$ws2.Range("C19").Value = 313

# --- Update sheet view selections ---
# Select on Hoja2 first (not the final active sheet) so it does not end up tabSelected
$ws2.Range("E16").Select()
# Activate Hoja1 (it was already the tabSelected sheet) and update its selection
$ws1.Activate()
$ws1.Range("F23").Select()

$wb.Save()
